# Apply updated cryptocurrency market data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Prefix with a single quote so Excel stores the exact string
    # (prevents auto-conversion of numeric-looking strings to numbers)
    $ws.Range($cellRef).Value = "'" + $text
}

Set-TextValue "D2" '42.262.04'
Set-TextValue "E2" '  -0.87%  '
Set-TextValue "D3" '2.242.83'
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "D5" '246.39'
Set-TextValue "E5" '  -1.64%  '
Set-TextValue "D6" '0.629'
Set-TextValue "E6" '  -1.78%  '
Set-TextValue "D7" '73.96'
Set-TextValue "E7" '  -3.02%  '
Set-TextValue "E8" '  +0.05%  '
Set-TextValue "D9" '0.616'
Set-TextValue "E9" '  -4.32%  '
Set-TextValue "D10" '42.02'
Set-TextValue "E10" '  +4.67%  '
Set-TextValue "D11" '0.0946'
Set-TextValue "E11" '  -2.67%  '
Set-TextValue "D12" '7.16'
Set-TextValue "E12" '  -2.15%  '
Set-TextValue "E13" '  -1.98%  '
Set-TextValue "D14" '14.46'
Set-TextValue "E14" '  -3.79%  '
Set-TextValue "E15" '  -1.42%  '
Set-TextValue "D16" '2.279.68'
Set-TextValue "E16" '  +0.54%  '
Set-TextValue "D17" '42.078.87'
Set-TextValue "E17" '  -1.07%  '
Set-TextValue "E18" '  -0.34%  '
Set-TextValue "E19" '  -0.53%  '
Set-TextValue "D20" '71.92'
Set-TextValue "E20" '  -0.20%  '
Set-TextValue "D21" '2.24'
Set-TextValue "E21" '  +4.01%  '
Set-TextValue "D22" '232.16'
Set-TextValue "E22" '  -1.52%  '
Set-TextValue "D23" '8.67'
Set-TextValue "E23" '  +32.58%  '
Set-TextValue "D25" '11.28'
Set-TextValue "E25" '  +0.07%  '
Set-TextValue "E26" '  -4.64%  '
Set-TextValue "E27" '  -3.33%  '
Set-TextValue "B28" 'Monero'
Set-TextValue "C28" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D28" '169.12'
Set-TextValue "E28" '  +1.11%  '
Set-TextValue "B29" 'Toncoin'
Set-TextValue "C29" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D29" '2.09'
Set-TextValue "E29" '  -2.29%  '
Set-TextValue "D31" '0.0823'
Set-TextValue "E31" '  -4.15%  '
Set-TextValue "E32" '  -4.06%  '
Set-TextValue "D33" '30.28'
Set-TextValue "E33" '  -3.96%  '
Set-TextValue "D34" '0.124'
Set-TextValue "E34" '  -2.37%  '
Set-TextValue "D35" '5.14'
Set-TextValue "E35" '  +8.60%  '
Set-TextValue "D36" '4.50'
Set-TextValue "E36" '  -0.49%  '
Set-TextValue "D37" '0.0305'
Set-TextValue "E37" '  -0.77%  '
Set-TextValue "D38" '13.68'
Set-TextValue "E38" '  -0.17%  '
Set-TextValue "E39" '  -3.40%  '
Set-TextValue "D40" '5.79'
Set-TextValue "E40" '  -1.12%  '
Set-TextValue "D41" '62.33'
Set-TextValue "E41" '  +1.11%  '
Set-TextValue "E42" '  -3.14%  '
Set-TextValue "D43" '107.61'
Set-TextValue "E43" '  +0.62%  '
Set-TextValue "E44" '  +1.90%  '
Set-TextValue "D45" '8.63'
Set-TextValue "E45" '  -2.50%  '
Set-TextValue "D46" '0.997'
Set-TextValue "E46" '  -0.22%  '
Set-TextValue "D47" '4.35'
Set-TextValue "E47" '  -8.26%  '
Set-TextValue "D48" '1.12'
Set-TextValue "E48" '  -2.72%  '
Set-TextValue "E49" '  -0.76%  '
Set-TextValue "E50" '  +1.27%  '
Set-TextValue "D51" '2.69'
Set-TextValue "E51" '  +0.08%  '
